$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Bird', ['Token Creature — Bird', 'Flying', '2/2'])"
$ws.Range("A3").Value = "('Boar', ['Token Creature — Boar', '2/2'])"
$ws.Range("A4").Value = "('Cleric', ['Token Enchantment Creature — Cleric', '2/1'])"
$ws.Range("A5").Value = "('Elemental', ['Token Creature — Elemental', '1/0'])"
$ws.Range("A6").Value = "(""Elspeth, Sun's Champion Emblem"", ['Emblem — Elspeth', 'Creatures you control get +2/+2 and have flying.'])"
$ws.Range("A7").Value = "('Golem', ['Token Enchantment Artifact Creature — Golem', '3/3'])"
$ws.Range("A8").Value = "('Harpy', ['Token Creature — Harpy', 'Flying', '1/1'])"
$ws.Range("A9").Value = "('Satyr', ['Token Creature — Satyr', '2/2'])"
$ws.Range("A10").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"

$ws.Range("A11:A36").Clear()
